$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.077.16"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.094.52"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'573.28"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "'176.75"
$ws.Range("E6").Value = "  +6.09%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.091.00"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "'6.50"
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("D11").Value = "'0.152"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").Value = "'0.466"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "'36.01"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "3.604.12"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "66.915.87"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "'7.01"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "3.088.24"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "'16.46"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "'484.17"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "'0.687"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'7.68"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'83.43"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Value = "'12.68"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D27").Value = "'10.33"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'7.85"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").Value = "'2.31"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "'2.59"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").Value = "'28.00"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "'0.112"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").Value = "0.0₃0933"
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("B36").Value = "Arweave"
$ws.Range("C36").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D36").Value = "'47.21"
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "'0.943"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "'5.56"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").Value = "'0.313"
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("D40").Value = "'49.05"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").Value = "'2.00"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "'8.25"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").Value = "'2.69"
$ws.Range("E44").Value = "  +10.02%  "
$ws.Range("D45").Value = "2.782.13"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "'368.68"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'135.68"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "'0.0342"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D50").Value = "'26.11"
$ws.Range("E50").Value = "  +7.42%  "
$ws.Range("D51").Value = "'2.29"
$ws.Range("E51").Value = "  +6.74%  "
